# Rewrite Sheet1 from the "Metadata" (single JSON object per row) sample
# into the "Line Detail" (JSON array of line items) sample, adding extra
# Company rows along the way.
#
# Target layout (A1:B6):
#   A1 Customer Name   | B1 Line Detail
#   A2 Company 1       | B2 [ {...Bolts...}, {...Smith...} ]
#   A3 Company 2       | B3 [ {...Braces...}, {...Wood...} ]
#   A4 Company 3       | B4 [{...Braces...}]
#   A5 Company 4       | (blank)
#   A6 Company 5       | B6 {...Braces...}   (single object, no array wrapper)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written row-by-row (A then B) so the shared-string table is
# built up in the same left-to-right, top-to-bottom order a human typing
# into the grid would produce.

# --- Row 1 (header) ---------------------------------------------------
$ws.Range("A1").Value = "Customer Name"
$ws.Range("B1").Value = "Line Detail"

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = "Company 1"
$ws.Range("B2").Value = '[ {"Id": "1", "Desc": "Bolts", "Amount": 101.15}, {"Id": "2", "Desc": "Smith", "Amount": 90.80} ]'

# --- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = "Company 2"
$ws.Range("B3").Value = '[ {"Id": "1", "Desc": "Braces", "Amount": 51.15}, {"Id": "2", "Desc": "Wood", "Amount": 190.10} ]'

# --- Row 4 --------------------------------------------------------------
$ws.Range("A4").Value = "Company 3"
$ws.Range("B4").Value = '[{"Id": "1", "Desc": "Braces", "Amount": 51.15}]'

# --- Row 5 (no Line Detail) ---------------------------------------------
$ws.Range("A5").Value = "Company 4"

# --- Row 6 ----------------------------------------------------------------
$ws.Range("A6").Value = "Company 5"
$ws.Range("B6").Value = '{"Id": "1", "Desc": "Braces", "Amount": 51.15}'

# --- Formatting --------------------------------------------------------
# B1 should carry the same bold header style as A1 (copy A1's format onto
# B1 so the engine reuses the existing bold cellXf instead of minting a
# new one).
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# B2 already carries the monospace "code" style from the source file;
# propagate that same format to the other JSON-bearing cells (and
# overwrite B3, which previously had an unused alignment-only style) so
# they all collapse onto the same cellXf.
$ws.Range("B2").Copy()
$ws.Range("B3:B4").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths -------------------------------------------------------
# Target widths (~50.16, ~12.66, ~11.83, ~14.16, 17 chars) land between this
# engine's sixth-of-a-character ColumnWidth quantization steps, so these
# inputs are chosen to round to the closest achievable step.
$ws.Columns.Item(1).ColumnWidth = 49.335
$ws.Columns.Item(2).ColumnWidth = 11.835
$ws.Columns.Item(3).ColumnWidth = 11.0025
$ws.Columns.Item(4).ColumnWidth = 13.335
$ws.Columns.Item(5).ColumnWidth = 16.1675

# --- Selection ---------------------------------------------------------
$ws.Range("B6").Select()

Write-Output "done"
